# Updated symbol list on Fri Dec 16 05:31:52 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) refresh -------------------------------------------------
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "261.36"
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "22.90"
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "6.207"
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.06145"
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "6.726"
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.467"
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "1.333"
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.7996"
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.1590"
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.08118"
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.03485"
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.03093"
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.09310"
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.862"
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.001694"
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.04788"
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.0006143"
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.006210"
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.004067"
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "3.693"
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.207"
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.1243"
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.0003204"
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.04616"
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.1118"
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.003131"
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.003323"
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.01023"
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.002971"
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.00005934"
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.7004"
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.08928"
$cell.Style = "Normal"

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.00002101"
$cell.Style = "Normal"

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.01011"
$cell.Style = "Normal"

# --- Column E (Volume(1h) label) updates --------------------------------------
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("E43").Value = "42KickTokenKICK"

# --- Rows 41/42/43: coin name + link rotate (Kick -> BKEX -> CEJI -> Kick) ----
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"

Write-Host "Done applying cryptos.xlsx symbol list update."
